$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 9
$ws.Range("H9").Value = 314.35715
$ws.Range("I9").Value = 344.21738
$ws.Range("J9").Value = 177
$ws.Range("K9").Value = 344.21738
$ws.Range("L9").Value = 177
$ws.Range("M9").Value = -175.21738
$ws.Range("N9").Value = -515
# row 32
$ws.Range("H32").Value = 477.3889
$ws.Range("I32").Value = 579.1
$ws.Range("J32").Value = 350.25
$ws.Range("K32").Value = 579.1
$ws.Range("L32").Value = 350.25
$ws.Range("M32").Value = -253.1
$ws.Range("N32").Value = -1002.25
# row 53
$ws.Range("H53").Value = 368.83334
$ws.Range("I53").Value = 562.9
$ws.Range("J53").Value = 126.25
$ws.Range("K53").Value = 562.9
$ws.Range("L53").Value = 126.25
$ws.Range("M53").Value = 74.10000000000002
$ws.Range("N53").Value = -1400.25
# row 125
$ws.Range("H125").Value = 1225.5
$ws.Range("I125").Value = 385
$ws.Range("J125").Value = 1393.6
$ws.Range("K125").Value = 3465
$ws.Range("L125").Value = 12542.4
$ws.Range("M125").Value = -1005
$ws.Range("N125").Value = -17462.4
# row 137
$ws.Range("H137").Value = 18376.807
$ws.Range("I137").Value = 657.2
$ws.Range("J137").Value = 60070
$ws.Range("K137").Value = 1971.6
$ws.Range("L137").Value = 180210
$ws.Range("M137").Value = 578.3999999999999
$ws.Range("N137").Value = -185310
# row 138
$ws.Range("H138").Value = 2633.27
$ws.Range("I138").Value = 1189.3
$ws.Range("J138").Value = 3595.9167
$ws.Range("K138").Value = 3567.9
$ws.Range("L138").Value = 10787.7501
$ws.Range("M138").Value = 1572.1
$ws.Range("N138").Value = -21067.7501

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 2
$ws.Range("H2").Value = 1201.625
$ws.Range("I2").Value = 1157.1428
$ws.Range("K2").Value = 1157.1428
$ws.Range("M2").Value = -1044.1428
# row 32
$ws.Range("H32").Value = 10189.12
$ws.Range("I32").Value = 7129.2324
$ws.Range("J32").Value = 28985.572
$ws.Range("K32").Value = 7129.2324
$ws.Range("L32").Value = 28985.572
$ws.Range("M32").Value = -6842.2324
$ws.Range("N32").Value = -29559.572
# row 116
$ws.Range("H116").Value = 1201.625
$ws.Range("I116").Value = 1157.1428
$ws.Range("K116").Value = 1157.1428
$ws.Range("M116").Value = 1136.8572

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 3
$ws.Range("H3").Value = 1201.625
$ws.Range("I3").Value = 1157.1428
$ws.Range("K3").Value = 1157.1428
$ws.Range("M3").Value = -1043.1428
# row 86
$ws.Range("H86").Value = 6252266.5
$ws.Range("I86").Value = 7694527.5
$ws.Range("J86").Value = 2469
$ws.Range("K86").Value = 7694527.5
$ws.Range("L86").Value = 2469
$ws.Range("M86").Value = -7693404.5
$ws.Range("N86").Value = -4715
# row 89
$ws.Range("H89").Value = 6252266.5
$ws.Range("I89").Value = 7694527.5
$ws.Range("J89").Value = 2469
$ws.Range("K89").Value = 38472637.5
$ws.Range("L89").Value = 12345
$ws.Range("M89").Value = -38467021.5
$ws.Range("N89").Value = -23577
# row 94
$ws.Range("H94").Value = 5174.1113
$ws.Range("I94").Value = 734.97144
$ws.Range("J94").Value = 20711.1
$ws.Range("K94").Value = 734.97144
$ws.Range("L94").Value = 20711.1
$ws.Range("M94").Value = -283.97144
$ws.Range("N94").Value = -21613.1
# row 105
$ws.Range("H105").Value = 2459.1538
$ws.Range("I105").Value = 2046.125
$ws.Range("J105").Value = 3120
$ws.Range("K105").Value = 2046.125
$ws.Range("L105").Value = 3120
$ws.Range("M105").Value = -299.125
$ws.Range("N105").Value = -6614
# row 141
$ws.Range("H141").Value = 34000
$ws.Range("J141").Value = 34000
$ws.Range("L141").Value = 34000
$ws.Range("N141").Value = -44360

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 16
$ws.Range("H16").Value = 2054.7856
$ws.Range("I16").Value = 1962.8889
$ws.Range("K16").Value = 1962.8889
$ws.Range("M16").Value = -1675.8889
# row 31
$ws.Range("H31").Value = 2149.6924
$ws.Range("I31").Value = 1231.4783
$ws.Range("K31").Value = 1231.4783
$ws.Range("M31").Value = -936.4783
# row 34
$ws.Range("H34").Value = 2149.6924
$ws.Range("I34").Value = 1231.4783
$ws.Range("K34").Value = 1231.4783
$ws.Range("M34").Value = -1029.4783
# row 113
$ws.Range("H113").Value = 2054.7856
$ws.Range("I113").Value = 1962.8889
$ws.Range("K113").Value = 1962.8889
$ws.Range("M113").Value = 207.1111000000001
# row 141
$ws.Range("H141").Value = 34096.727
$ws.Range("J141").Value = 35244.19
$ws.Range("L141").Value = 35244.19
$ws.Range("N141").Value = -45604.19

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 90018.87
$ws.Range("J2").Value = 23.916666
$ws.Range("L2").Value = 143.499996
$ws.Range("N2").Value = -369.499996
# row 5
$ws.Range("H5").Value = 857.17145
$ws.Range("J5").Value = 1311
$ws.Range("L5").Value = 3933
$ws.Range("N5").Value = -4157
# row 23
$ws.Range("H23").Value = 49.846153
$ws.Range("I23").Value = 9
$ws.Range("J23").Value = 57.272728
$ws.Range("K23").Value = 27
$ws.Range("L23").Value = 171.818184
$ws.Range("M23").Value = 208
$ws.Range("N23").Value = -641.818184
# row 34
$ws.Range("H34").Value = 489.18182
$ws.Range("I34").Value = 50.923077
$ws.Range("J34").Value = 1122.2222
$ws.Range("K34").Value = 152.769231
$ws.Range("L34").Value = 3366.6666
$ws.Range("M34").Value = -68.76923099999999
$ws.Range("N34").Value = -3534.6666
# row 40
$ws.Range("H40").Value = 5259.421
$ws.Range("I40").Value = 52.07143
$ws.Range("J40").Value = 19840
$ws.Range("K40").Value = 208.28572
$ws.Range("L40").Value = 79360
$ws.Range("M40").Value = -139.28572
$ws.Range("N40").Value = -79498
# row 113
$ws.Range("H113").Value = 1301.5
$ws.Range("I113").Value = 348.875
$ws.Range("J113").Value = 2148.2778
$ws.Range("K113").Value = 1046.625
$ws.Range("L113").Value = 6444.8334
$ws.Range("M113").Value = 1123.375
$ws.Range("N113").Value = -10784.8334
# row 135
$ws.Range("H135").Value = 857.17145
$ws.Range("J135").Value = 1311
$ws.Range("L135").Value = 11799
$ws.Range("N135").Value = -16869

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 5
$ws.Range("H5").Value = 10000
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = $null
# row 80
$ws.Range("H80").Value = 3000
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -4996
# row 83
$ws.Range("H83").Value = 3000
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -24984
# row 113
$ws.Range("H113").Value = 1479.8334
$ws.Range("I113").Value = 960
$ws.Range("J113").Value = 1999.6666
$ws.Range("K113").Value = 960
$ws.Range("L113").Value = 1999.6666
$ws.Range("M113").Value = 1210
$ws.Range("N113").Value = -6339.6666

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 2260.9092
$ws.Range("I7").Value = 1746.25
$ws.Range("J7").Value = 3633.3333
$ws.Range("K7").Value = 1746.25
$ws.Range("L7").Value = 3633.3333
$ws.Range("M7").Value = -1634.25
$ws.Range("N7").Value = -3857.3333
# row 9
$ws.Range("H9").Value = 1996.125
$ws.Range("I9").Value = 193.8
$ws.Range("J9").Value = 5000
$ws.Range("K9").Value = 193.8
$ws.Range("L9").Value = 5000
$ws.Range("M9").Value = 30.19999999999999
$ws.Range("N9").Value = -5448
# row 16
$ws.Range("H16").Value = 976.4286
$ws.Range("I16").Value = 1121.25
$ws.Range("J16").Value = 783.3333
$ws.Range("K16").Value = 1121.25
$ws.Range("L16").Value = 783.3333
$ws.Range("M16").Value = -951.25
$ws.Range("N16").Value = -1123.3333
# row 35
$ws.Range("H35").Value = 2213
$ws.Range("I35").Value = 257.5
$ws.Range("K35").Value = 257.5
$ws.Range("M35").Value = 78.5
# row 40
$ws.Range("H40").Value = 1714.9375
$ws.Range("I40").Value = 1730.3636
$ws.Range("K40").Value = 1730.3636
$ws.Range("M40").Value = -1594.3636
# row 126
$ws.Range("H126").Value = 2260.9092
$ws.Range("I126").Value = 1746.25
$ws.Range("J126").Value = 3633.3333
$ws.Range("K126").Value = 5238.75
$ws.Range("L126").Value = 10899.9999
$ws.Range("M126").Value = -2768.75
$ws.Range("N126").Value = -15839.9999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 11
$ws.Range("H11").Value = 500000
$ws.Range("I11").Value = 500000
$ws.Range("K11").Value = 500000
$ws.Range("M11").Value = -499858
# row 51
$ws.Range("H51").Value = 121207.78
$ws.Range("I51").Value = 253217.5
$ws.Range("J51").Value = 15600
$ws.Range("K51").Value = 253217.5
$ws.Range("L51").Value = 15600
$ws.Range("M51").Value = -252707.5
$ws.Range("N51").Value = -16620
# row 100
$ws.Range("H100").Value = 4902.2915
$ws.Range("I100").Value = 10145.5
$ws.Range("J100").Value = 1157.1428
$ws.Range("K100").Value = 20291
$ws.Range("L100").Value = 2314.2856
$ws.Range("M100").Value = -19750
$ws.Range("N100").Value = -3396.2856
